# Auto-generated edit script: update leve-profit figures across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3540
$ws.Range("I64").Value = 3320
$ws.Range("J64").Value = 3650
$ws.Range("K64").Value = 3320
$ws.Range("L64").Value = 3650
$ws.Range("M64").Value = -3072
$ws.Range("N64").Value = -4146

# Row 67
$ws.Range("H67").Value = 3540
$ws.Range("I67").Value = 3320
$ws.Range("J67").Value = 3650
$ws.Range("K67").Value = 3320
$ws.Range("L67").Value = 3650
$ws.Range("M67").Value = -2462
$ws.Range("N67").Value = -5366

# Row 138
$ws.Range("H138").Value = 2603.1082
$ws.Range("I138").Value = 2348.6667
$ws.Range("J138").Value = 2725.24
$ws.Range("K138").Value = 7046.000100000001
$ws.Range("L138").Value = 8175.719999999999
$ws.Range("M138").Value = -1906.000100000001
$ws.Range("N138").Value = -18455.72

# Row 141
$ws.Range("H141").Value = 3175.3333
$ws.Range("I141").Value = 1794.1
$ws.Range("J141").Value = 4901.875
$ws.Range("K141").Value = 5382.299999999999
$ws.Range("L141").Value = 14705.625
$ws.Range("M141").Value = -202.2999999999993
$ws.Range("N141").Value = -25065.625

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1765.3334
$ws.Range("I45").Value = 1878.5
$ws.Range("J45").Value = 1674.8
$ws.Range("K45").Value = 1878.5
$ws.Range("L45").Value = 1674.8
$ws.Range("M45").Value = -1501.5
$ws.Range("N45").Value = -2428.8

# Row 74
$ws.Range("H74").Value = 1640.58
$ws.Range("I74").Value = 1288.2821
$ws.Range("J74").Value = 2889.6365
$ws.Range("K74").Value = 1288.2821
$ws.Range("L74").Value = 2889.6365
$ws.Range("M74").Value = -414.2820999999999

# Row 77
$ws.Range("H77").Value = 1640.58
$ws.Range("I77").Value = 1288.2821
$ws.Range("J77").Value = 2889.6365
$ws.Range("K77").Value = 6441.4105
$ws.Range("L77").Value = 14448.1825
$ws.Range("M77").Value = -2073.4105

# Row 132
$ws.Range("H132").Value = 33335774
$ws.Range("I132").Value = 45456056
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 136368168
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -136365638
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
# Row 55
$ws.Range("H55").Value = 20000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 20000
$ws.Range("N55").Value = -20546

# Row 64
$ws.Range("H64").Value = 654.8
$ws.Range("I64").Value = 847
$ws.Range("J64").Value = 462.6
$ws.Range("K64").Value = 847
$ws.Range("L64").Value = 462.6
$ws.Range("M64").Value = -622
$ws.Range("N64").Value = -912.6

# Row 67
$ws.Range("H67").Value = 654.8
$ws.Range("I67").Value = 847
$ws.Range("J67").Value = 462.6
$ws.Range("K67").Value = 847
$ws.Range("L67").Value = 462.6
$ws.Range("M67").Value = -67
$ws.Range("N67").Value = -2022.6

# Row 86
$ws.Range("H86").Value = 2471.2856
$ws.Range("I86").Value = 2399.8
$ws.Range("J86").Value = 2650
$ws.Range("K86").Value = 2399.8
$ws.Range("L86").Value = 2650
$ws.Range("M86").Value = -1276.8
$ws.Range("N86").Value = -4896

# Row 89
$ws.Range("H89").Value = 2471.2856
$ws.Range("I89").Value = 2399.8
$ws.Range("J89").Value = 2650
$ws.Range("K89").Value = 11999
$ws.Range("L89").Value = 13250
$ws.Range("M89").Value = -6383
$ws.Range("N89").Value = -24482

# Row 107
$ws.Range("H107").Value = 1757.0526
$ws.Range("I107").Value = 1544.4166
$ws.Range("J107").Value = 2121.5715
$ws.Range("K107").Value = 1544.4166
$ws.Range("L107").Value = 2121.5715
$ws.Range("M107").Value = 375.5834
$ws.Range("N107").Value = -5961.5715

# Row 134
$ws.Range("H134").Value = 2356.5754
$ws.Range("I134").Value = 1430.0698
$ws.Range("J134").Value = 3684.5667
$ws.Range("K134").Value = 4290.2094
$ws.Range("L134").Value = 11053.7001
$ws.Range("M134").Value = -1755.2094

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2800.5833
$ws.Range("I62").Value = 2600.125
$ws.Range("J62").Value = 3201.5
$ws.Range("K62").Value = 2600.125
$ws.Range("L62").Value = 3201.5
$ws.Range("M62").Value = -1976.125

# Row 65
$ws.Range("H65").Value = 2800.5833
$ws.Range("I65").Value = 2600.125
$ws.Range("J65").Value = 3201.5
$ws.Range("K65").Value = 13000.625
$ws.Range("L65").Value = 16007.5
$ws.Range("M65").Value = -9880.625

# Row 98
$ws.Range("H98").Value = 20500
$ws.Range("I98").Value = 22000
$ws.Range("J98").Value = 20384.615
$ws.Range("K98").Value = 22000
$ws.Range("L98").Value = 20384.615
$ws.Range("M98").Value = -19754
$ws.Range("N98").Value = -24876.615

# Row 134
$ws.Range("H134").Value = 1722326.9
$ws.Range("I134").Value = 2818186.2
$ws.Range("J134").Value = 352502.75
$ws.Range("K134").Value = 8454558.600000001
$ws.Range("L134").Value = 1057508.25
$ws.Range("M134").Value = -8452023.600000001
$ws.Range("N134").Value = -1062578.25

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 355.94446
$ws.Range("I107").Value = 366.64285
$ws.Range("J107").Value = 318.5
$ws.Range("K107").Value = 1099.92855
$ws.Range("L107").Value = 955.5
$ws.Range("M107").Value = 820.0714499999999
$ws.Range("N107").Value = -4795.5

# Row 131
$ws.Range("H131").Value = 875.52
$ws.Range("I131").Value = 630
$ws.Range("J131").Value = 878
$ws.Range("K131").Value = 1890
$ws.Range("L131").Value = 2634
$ws.Range("M131").Value = 3150
$ws.Range("N131").Value = -12714

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1844.5555
$ws.Range("I113").Value = 2022.2
$ws.Range("J113").Value = 1622.5
$ws.Range("K113").Value = 2022.2
$ws.Range("L113").Value = 1622.5
$ws.Range("M113").Value = 147.8
$ws.Range("N113").Value = -5962.5

# Row 132
$ws.Range("H132").Value = 33337948
$ws.Range("I132").Value = 62504890
$ws.Range("J132").Value = 4300.0713
$ws.Range("K132").Value = 187514670
$ws.Range("L132").Value = 12900.2139
$ws.Range("M132").Value = -187512140
$ws.Range("N132").Value = -17960.2139

# Row 140
$ws.Range("H140").Value = 40000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 40000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2733.611
$ws.Range("I61").Value = 2875.2222
$ws.Range("J61").Value = 2592
$ws.Range("K61").Value = 2875.2222
$ws.Range("L61").Value = 2592
$ws.Range("M61").Value = -2673.2222
$ws.Range("N61").Value = -2996

# Row 113
$ws.Range("H113").Value = 2733.611
$ws.Range("I113").Value = 2875.2222
$ws.Range("J113").Value = 2592
$ws.Range("K113").Value = 2875.2222
$ws.Range("L113").Value = 2592
$ws.Range("M113").Value = -705.2222000000002
$ws.Range("N113").Value = -6932

# Row 132
$ws.Range("H132").Value = 3096.1292
$ws.Range("I132").Value = 2166.2856
$ws.Range("J132").Value = 5048.8
$ws.Range("K132").Value = 6498.8568
$ws.Range("L132").Value = 15146.4
$ws.Range("M132").Value = -3968.8568
$ws.Range("N132").Value = -20206.4

# Row 136
$ws.Range("H136").Value = 1961.9656
$ws.Range("I136").Value = 1427.3636
$ws.Range("J136").Value = 3642.1428
$ws.Range("K136").Value = 4282.0908
$ws.Range("L136").Value = 10926.4284
$ws.Range("M136").Value = -1732.0908
$ws.Range("N136").Value = -16026.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 18138.8
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 18138.8
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 18138.8
$ws.Range("N70").Value = -18768.8
$ws.Range("M70").ClearContents()

# Row 73
$ws.Range("H73").Value = 18138.8
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 18138.8
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 18138.8
$ws.Range("N73").Value = -20322.8
$ws.Range("M73").ClearContents()

# Row 122
$ws.Range("H122").Value = 1191145
$ws.Range("I122").Value = 1681230.6
$ws.Range("J122").Value = 937.1429000000001
$ws.Range("K122").Value = 5043691.800000001
$ws.Range("L122").Value = 2811.4287
$ws.Range("M122").Value = -5041241.800000001
